$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = '59.040.59'
$ws.Cells.Item(2,5).Value = '  +1.38%  '

# Row 3
$ws.Cells.Item(3,4).Value = '2.586.98'
$ws.Cells.Item(3,5).Value = '  -1.07%  '

# Row 4
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = '1.00'
$ws.Cells.Item(4,4).NumberFormat = "General"
$ws.Cells.Item(4,5).Value = '  +0.21%  '

# Row 5
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '522.68'
$ws.Cells.Item(5,4).NumberFormat = "General"
$ws.Cells.Item(5,5).Value = '  -0.21%  '

# Row 6
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '139.31'
$ws.Cells.Item(6,4).NumberFormat = "General"
$ws.Cells.Item(6,5).Value = '  -3.62%  '

# Row 7
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = '0.999'
$ws.Cells.Item(7,4).NumberFormat = "General"

# Row 9
$ws.Cells.Item(9,4).Value = '2.599.56'
$ws.Cells.Item(9,5).Value = '  -0.51%  '

# Row 10
$ws.Cells.Item(10,5).Value = '  -1.70%  '

# Row 11
$ws.Cells.Item(11,5).Value = '  -0.54%  '

# Row 12
$ws.Cells.Item(12,5).Value = '  -1.37%  '

# Row 13
$ws.Cells.Item(13,5).Value = '  +3.02%  '

# Row 14
$ws.Cells.Item(14,4).Value = '3.050.48'
$ws.Cells.Item(14,5).Value = '  -0.12%  '

# Row 15
$ws.Cells.Item(15,4).Value = '58.980.29'
$ws.Cells.Item(15,5).Value = '  +1.36%  '

# Row 16
$ws.Cells.Item(16,5).Value = '  +0.06%  '

# Row 17
$ws.Cells.Item(17,4).Value = '2.614.78'
$ws.Cells.Item(17,5).Value = '  +0.51%  '

# Row 18
$ws.Cells.Item(18,5).Value = '  -1.09%  '

# Row 19
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '340.16'
$ws.Cells.Item(19,4).NumberFormat = "General"
$ws.Cells.Item(19,5).Value = '  +0.03%  '

# Row 20
$ws.Cells.Item(20,5).Value = '  -1.45%  '

# Row 21
$ws.Cells.Item(21,5).Value = '  -2.39%  '

# Row 22
$ws.Cells.Item(22,5).Value = '  +0.71%  '

# Row 23
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '1.00'
$ws.Cells.Item(23,4).NumberFormat = "General"
$ws.Cells.Item(23,5).Value = '  +0.29%  '

# Row 24
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '66.36'
$ws.Cells.Item(24,4).NumberFormat = "General"
$ws.Cells.Item(24,5).Value = '  +1.27%  '

# Row 25
$ws.Cells.Item(25,5).Value = '  +1.01%  '

# Row 26
$ws.Cells.Item(26,5).Value = '  -0.11%  '

# Row 27
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '0.999'
$ws.Cells.Item(27,4).NumberFormat = "General"
$ws.Cells.Item(27,5).Value = '  +0.28%  '

# Row 28
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = '7.06'
$ws.Cells.Item(28,4).NumberFormat = "General"
$ws.Cells.Item(28,5).Value = '  +0.20%  '

# Row 29
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '0.999'
$ws.Cells.Item(29,4).NumberFormat = "General"
$ws.Cells.Item(29,5).Value = '  +0.07%  '

# Row 30
$ws.Cells.Item(30,5).Value = '  -3.83%  '

# Row 31
$ws.Cells.Item(31,5).Value = '  -6.09%  '

# Row 32
$ws.Cells.Item(32,5).Value = '  -0.54%  '

# Row 33
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '18.71'
$ws.Cells.Item(33,4).NumberFormat = "General"
$ws.Cells.Item(33,5).Value = '  -0.92%  '

# Row 34
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = '149.33'
$ws.Cells.Item(34,4).NumberFormat = "General"
$ws.Cells.Item(34,5).Value = '  -0.33%  '

# Row 35
$ws.Cells.Item(35,5).Value = '  -2.37%  '

# Row 36
$ws.Cells.Item(36,5).Value = '  -2.39%  '

# Row 37
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '36.75'
$ws.Cells.Item(37,4).NumberFormat = "General"
$ws.Cells.Item(37,5).Value = '  +1.72%  '

# Row 38
$ws.Cells.Item(38,5).Value = '  +0.25%  '

# Row 39
$ws.Cells.Item(39,5).Value = '  -3.10%  '

# Row 40
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '0.811'
$ws.Cells.Item(40,4).NumberFormat = "General"
$ws.Cells.Item(40,5).Value = '  -6.41%  '

# Row 41
$ws.Cells.Item(41,5).Value = '  -1.13%  '

# Row 42
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '0.998'
$ws.Cells.Item(42,4).NumberFormat = "General"
$ws.Cells.Item(42,5).Value = '  -0.02%  '

# Row 43
$ws.Cells.Item(43,2).Value = 'Mantle'
$ws.Cells.Item(43,3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '0.603'
$ws.Cells.Item(43,4).NumberFormat = "General"
$ws.Cells.Item(43,5).Value = '  +0.94%  '

# Row 44
$ws.Cells.Item(44,2).Value = 'Bittensor'
$ws.Cells.Item(44,3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '271.97'
$ws.Cells.Item(44,4).NumberFormat = "General"
$ws.Cells.Item(44,5).Value = '  -0.70%  '

# Row 45
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '10.77'
$ws.Cells.Item(45,4).NumberFormat = "General"
$ws.Cells.Item(45,5).Value = '  +1.12%  '

# Row 46
$ws.Cells.Item(46,5).Value = '  -0.87%  '

# Row 47
$ws.Cells.Item(47,5).Value = '  -1.54%  '

# Row 48
$ws.Cells.Item(48,5).Value = '  -2.41%  '

# Row 49
$ws.Cells.Item(49,4).Value = '1.968.91'
$ws.Cells.Item(49,5).Value = '  -0.66%  '

# Row 50
$ws.Cells.Item(50,2).Value = 'VeChain'
$ws.Cells.Item(50,3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = '0.0221'
$ws.Cells.Item(50,4).NumberFormat = "General"
$ws.Cells.Item(50,5).Value = '  -0.86%  '

# Row 51
$ws.Cells.Item(51,2).Value = 'RenderToken'
$ws.Cells.Item(51,3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = '4.48'
$ws.Cells.Item(51,4).NumberFormat = "General"
$ws.Cells.Item(51,5).Value = '  -4.03%  '
